$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Fill in "Hand in" / "Test AUC" for the two rows that already existed
# (row 19 -> "2023-03-03-2029_RF_C1.csv", row 20 -> "2023-03-03-2029_RF_C5.csv")
$ws.Range("M19").Value = "March 5, 2023, 8:36 p.m."
$ws.Range("N19").Value = 0.518

$ws.Range("M20").Value = "March 5, 2023, 8:51 p.m."
$ws.Range("N20").Value = 0.666

# Add the two new submissions as new table rows
$newRow1 = $lo.ListRows.Add()
$newRow2 = $lo.ListRows.Add()

# Row 23 - tuned Logistic Regression, center-only submission
$ws.Range("M23").Value = "March 7, 2023, 9:56 p.m."
$ws.Range("H23").Value = "{'C': 0.6, 'max_iter': 100, 'penalty': 'l2', 'solver': 'lbfgs'}"
$ws.Range("A23").Value = "2023-03-07-2255_tuned_LogReg.csv"
$ws.Range("B23").Value = "Logistic Classifier"
$ws.Range("C23").Value = "MoCo"
$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "weakly supervision"
$ws.Range("N23").Value = 0.646

# Row 24 - Logistic Regression C5, center-only submission
$ws.Range("M24").Value = "March 7, 2023, 5:46 p.m."
$ws.Range("A24").Value = "2023-03-07-1846_LogReg_C5.csv"
$ws.Range("B24").Value = "Logistic Classifier"
$ws.Range("C24").Value = "MoCo"
$ws.Range("D24").Value = "-"
$ws.Range("E24").Value = "-"
$ws.Range("F24").Value = "weakly supervision"
$ws.Range("N24").Value = 0.624

# Reset the view: clear the frozen/scrolled left column and move the
# selection to the first empty row below the table
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A25").Select()
